# HERCULES-8990 - [UPV] Multi-idioma - Carga de traducciones - Etiquetas y mensajes
$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# 1) document.xml: shrink/grow a couple of table column widths (dxa -> pt/20)
# ---------------------------------------------------------------------------
$tbl2 = $d.Tables.Item(2)
$tbl2.Rows.Item(2).Cells.Item(1).Width = 240.6   # 4814 -> 4812 dxa
$tbl2.Rows.Item(2).Cells.Item(3).Width = 116.4   # 2326 -> 2328 dxa

$tbl3 = $d.Tables.Item(3)
$tbl3.Rows.Item(1).Cells.Item(1).Width = 240.6   # 4814 -> 4812 dxa
$tbl3.Rows.Item(1).Cells.Item(3).Width = 116.4   # 2326 -> 2328 dxa

Write-Host "table widths updated"

# ---------------------------------------------------------------------------
# 2) document.xml: sectPr page margins (top/bottom grow slightly)
# ---------------------------------------------------------------------------
$ps = $d.Sections.Item(1).PageSetup
$ps.TopMargin = 113.4      # 1977 -> 2268 dxa
$ps.BottomMargin = 100.95  # 1969 -> 2019 dxa

Write-Host "page margins updated"
